$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New arrival rows for "Sunday, Jan 08" (data downloaded from 11 airports)
# Columns: A=NUMBER B=DATE C=TIME D=FLIGHT E=FROM F=SHORT G=AIRLINE H=MODEL I=AIRCFAT ID J=STATUS K=(blank) L=DIFFERENCE M=(blank)

$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(35, 3).Value = "9:50 PM"
$ws.Cells.Item(35, 4).Value = "FR7678"
$ws.Cells.Item(35, 5).Value = "Stockholm"
$ws.Cells.Item(35, 6).Value = "(ARN)"
$ws.Cells.Item(35, 7).Value = "Ryanair "
$ws.Cells.Item(35, 8).Value = "B38M"
$ws.Cells.Item(35, 9).Value = "(9H-VUJ)"
$ws.Cells.Item(35, 10).Value = "9:42 PM"
$ws.Cells.Item(35, 12).Value = "0 hours, -8 minutes"

$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(36, 3).Value = "10:05 PM"
$ws.Cells.Item(36, 4).Value = "LH1380"
$ws.Cells.Item(36, 5).Value = "Frankfurt"
$ws.Cells.Item(36, 6).Value = "(FRA)"
$ws.Cells.Item(36, 7).Value = "Lufthansa "
$ws.Cells.Item(36, 8).Value = "CRJ9"
$ws.Cells.Item(36, 9).Value = "(D-ACNP)"
$ws.Cells.Item(36, 10).Value = "9:53 PM"
$ws.Cells.Item(36, 12).Value = "0 hours, -12 minutes"

$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(37, 3).Value = "10:10 PM"
$ws.Cells.Item(37, 4).Value = "FR7945"
$ws.Cells.Item(37, 5).Value = "Leeds"
$ws.Cells.Item(37, 6).Value = "(LBA)"
$ws.Cells.Item(37, 7).Value = "Ryanair "
$ws.Cells.Item(37, 8).Value = "B738"
$ws.Cells.Item(37, 9).Value = "(SP-RKR)"
$ws.Cells.Item(37, 10).Value = "9:40 PM"
$ws.Cells.Item(37, 12).Value = "0 hours, -30 minutes"

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(38, 3).Value = "10:45 PM"
$ws.Cells.Item(38, 4).Value = "FR6264"
$ws.Cells.Item(38, 5).Value = "Krakow"
$ws.Cells.Item(38, 6).Value = "(KRK)"
$ws.Cells.Item(38, 7).Value = "Ryanair "
$ws.Cells.Item(38, 8).Value = "B738"
$ws.Cells.Item(38, 9).Value = "(SP-RSM)"
$ws.Cells.Item(38, 10).Value = "10:24 PM"
$ws.Cells.Item(38, 12).Value = "0 hours, -21 minutes"
